$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = 6.5
$ws.Range("K2").Value = 2.38
$ws.Range("L2").Value = 6.5
$ws.Range("M2").Value = 1.04
$ws.Range("N2").Value = 12
$ws.Range("O2").Value = 1.25
$ws.Range("P2").Value = 4
$ws.Range("Q2").Value = 1.8
$ws.Range("R2").Value = 2
$ws.Range("S2").Value = 3
$ws.Range("T2").Value = 1.4
$ws.Range("U2").Value = 1.33
$ws.Range("V2").Value = 3.25
$ws.Range("W2").Value = 1.95
$ws.Range("X2").Value = 1.8
$ws.Range("Y2").Value = 7
$ws.Range("Z2").Value = 7
$ws.Range("AD2").Value = 26
$ws.Range("AE2").Value = 12
$ws.Range("AG2").Value = 19
$ws.Range("AI2").Value = 351
$ws.Range("AJ2").Value = 15
$ws.Range("AN2").Value = 51
$ws.Range("AO2").Value = 51

# Row 3
$ws.Range("G3").Value = 2.25
$ws.Range("I3").Value = 3.5
$ws.Range("J3").Value = 3.1
$ws.Range("K3").Value = 1.95
$ws.Range("W3").Value = 2.05
$ws.Range("X3").Value = 1.7
$ws.Range("AK3").Value = 15
$ws.Range("AP3").Value = 1.98
$ws.Range("AQ3").Value = 1.92

# Row 5
$ws.Range("G5").Value = 2.1
$ws.Range("H5").Value = 3
$ws.Range("I5").Value = 3.8
$ws.Range("L5").Value = 4.75
$ws.Range("Z5").Value = 8.5
$ws.Range("AB5").Value = 19
$ws.Range("AJ5").Value = 8
$ws.Range("AK5").Value = 17
$ws.Range("AP5").Value = 2.1
$ws.Range("AQ5").Value = 1.78

# Row 9
$ws.Range("G9").Value = 2.63
$ws.Range("H9").Value = 3.25
$ws.Range("I9").Value = 2.55
$ws.Range("J9").Value = 3.5
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 3.4
$ws.Range("M9").Value = 1.07
$ws.Range("N9").Value = 9
$ws.Range("O9").Value = 1.4
$ws.Range("P9").Value = 2.75
$ws.Range("Q9").Value = 2.3
$ws.Range("R9").Value = 1.6
$ws.Range("S9").Value = 4.33
$ws.Range("T9").Value = 1.2
$ws.Range("U9").Value = 1.5
$ws.Range("V9").Value = 2.5
$ws.Range("Y9").Value = 7.5
$ws.Range("Z9").Value = 12
$ws.Range("AB9").Value = 26
$ws.Range("AC9").Value = 23
$ws.Range("AE9").Value = 8
$ws.Range("AK9").Value = 12
$ws.Range("AM9").Value = 26
$ws.Range("AN9").Value = 23

# Row 15
$ws.Range("G15").Value = 1.38
$ws.Range("I15").Value = 7
$ws.Range("L15").Value = 6.5
$ws.Range("Q15").Value = 1.53
$ws.Range("R15").Value = 2.4
$ws.Range("S15").Value = 2.25
$ws.Range("T15").Value = 1.57
$ws.Range("Y15").Value = 8.5
$ws.Range("Z15").Value = 7.5
$ws.Range("AB15").Value = 9.5
$ws.Range("AO15").Value = 41
$ws.Range("AR15").Value = 1.9
$ws.Range("AS15").Value = 1.95

# Row 17
$ws.Range("G17").Value = 3.3
$ws.Range("H17").Value = 3.9
$ws.Range("I17").Value = 2
$ws.Range("J17").Value = 3.6
$ws.Range("K17").Value = 2.4
$ws.Range("L17").Value = 2.6
$ws.Range("M17").Value = 1.02
$ws.Range("N17").Value = 19
$ws.Range("O17").Value = 1.14
$ws.Range("P17").Value = 5.5
$ws.Range("Q17").Value = 1.53
$ws.Range("R17").Value = 2.4
$ws.Range("S17").Value = 2.25
$ws.Range("T17").Value = 1.57
$ws.Range("U17").Value = 1.29
$ws.Range("V17").Value = 3.5
$ws.Range("W17").Value = 1.5
$ws.Range("X17").Value = 2.5
$ws.Range("Y17").Value = 15
$ws.Range("Z17").Value = 21
$ws.Range("AA17").Value = 12
$ws.Range("AB17").Value = 34
$ws.Range("AC17").Value = 23
$ws.Range("AD17").Value = 26
$ws.Range("AE17").Value = 19
$ws.Range("AF17").Value = 8
$ws.Range("AH17").Value = 34
$ws.Range("AI17").Value = 101
$ws.Range("AJ17").Value = 11
$ws.Range("AK17").Value = 12
$ws.Range("AM17").Value = 19
$ws.Range("AN17").Value = 15
$ws.Range("AO17").Value = 19

# Row 18
$ws.Range("G18").Value = 6.5
$ws.Range("I18").Value = 1.36
$ws.Range("J18").Value = 6
$ws.Range("L18").Value = 1.8
$ws.Range("N18").Value = 23
$ws.Range("W18").Value = 1.53
$ws.Range("X18").Value = 2.38
$ws.Range("AF18").Value = 11
$ws.Range("AG18").Value = 15
$ws.Range("AI18").Value = 101
$ws.Range("AJ18").Value = 13

# Row 19
$ws.Range("G19").Value = 2.7
$ws.Range("I19").Value = 2.4
$ws.Range("J19").Value = 3.1
$ws.Range("L19").Value = 2.88
$ws.Range("O19").Value = 1.17
$ws.Range("P19").Value = 5
$ws.Range("S19").Value = 2.25
$ws.Range("T19").Value = 1.57
$ws.Range("Z19").Value = 17
$ws.Range("AC19").Value = 19
$ws.Range("AD19").Value = 21
$ws.Range("AJ19").Value = 13
$ws.Range("AL19").Value = 10
$ws.Range("AN19").Value = 17

# Row 20
$ws.Range("G20").Value = 2.75
$ws.Range("I20").Value = 2.63
$ws.Range("J20").Value = 3.25
$ws.Range("L20").Value = 3.2
$ws.Range("O20").Value = 1.29
$ws.Range("P20").Value = 3.5
$ws.Range("Q20").Value = 1.98
$ws.Range("R20").Value = 1.88
$ws.Range("Y20").Value = 9.5
$ws.Range("AA20").Value = 11
$ws.Range("AB20").Value = 26
$ws.Range("AG20").Value = 12
$ws.Range("AJ20").Value = 9.5
$ws.Range("AK20").Value = 13
$ws.Range("AL20").Value = 10
$ws.Range("AM20").Value = 26
$ws.Range("AP20").Value = 1.46
$ws.Range("AQ20").Value = 2.75

# Row 22
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = 3.3
$ws.Range("I22").Value = 3.7
$ws.Range("K22").Value = 2.1
$ws.Range("N22").Value = 9.5
$ws.Range("O22").Value = 1.3
$ws.Range("P22").Value = 3.4
$ws.Range("Q22").Value = 2.05
$ws.Range("R22").Value = 1.8
$ws.Range("S22").Value = 3.5
$ws.Range("T22").Value = 1.29
$ws.Range("W22").Value = 1.8
$ws.Range("X22").Value = 1.95
$ws.Range("Y22").Value = 7.5
$ws.Range("AB22").Value = 17
$ws.Range("AE22").Value = 9.5
$ws.Range("AF22").Value = 6.5
$ws.Range("AI22").Value = 251
$ws.Range("AJ22").Value = 11
$ws.Range("AN22").Value = 29

# Row 24
$ws.Range("G24").Value = 1.65
$ws.Range("H24").Value = 3.75
$ws.Range("I24").Value = 5.25
$ws.Range("J24").Value = 2.3
$ws.Range("L24").Value = 6.5
$ws.Range("N24").Value = 8
$ws.Range("Y24").Value = 5
$ws.Range("Z24").Value = 6.5
$ws.Range("AJ24").Value = 11
$ws.Range("AK24").Value = 26
$ws.Range("AL24").Value = 19
$ws.Range("AM24").Value = 67
$ws.Range("AN24").Value = 51
